$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently carries two header rows (row 1 = partial/merged-cell
# leftovers, row 2 = unit labels) above 10 data rows (rows 3-12). Deleting the
# old units row merges everything into a single header row and shifts the
# data rows up by one (rows 3-12 -> rows 2-11).
$ws.Rows("2:2").Delete() | Out-Null

# --- Row 1: rewrite as a single header row ---
# Columns A-E get new plain-text headers with default (General) formatting.
$hdrLeft = $ws.Range("A1:E1")
$hdrLeft.ClearFormats() | Out-Null
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

# Columns F-K get the unit headers, styled with the Arial 9 font used
# throughout the data rows (no explicit number format override).
$hdrRight = $ws.Range("F1:K1")
$hdrRight.ClearFormats() | Out-Null
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

$hdrStyle = $wb.Styles.Add("HdrStyle")
$hdrStyle.Font.Name = "Arial"
$hdrStyle.Font.Size = 9
$hdrRight.Style = "HdrStyle"
$wb.Styles.Item("HdrStyle").Delete() | Out-Null

# --- Update the selection shown when the workbook is (re)opened ---
$ws.Range("A2:K2").Select() | Out-Null

Write-Output "edit complete"
